# Updated symbol list on Fri Jan 27 22:25:27 UTC 2023 with GitHub Actions
#
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# crypto rows. The source sheet stores these as literal text (e.g.
# "307.82", "0.93%"), so a plain Range.Value assignment would let Excel's
# auto-detection convert them to numbers/percentages and reformat the
# cell. To keep every updated cell as plain text -- identical in kind and
# style to how it was authored -- we stage each new value in a scratch
# cell with a leading apostrophe (forces text), copy it, and
# Paste-Special "Values only" onto the real target. That brings over just
# the literal text, with no number/percent formatting riding along.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Addr = "D2"; Val = "307.82" },
    @{ Addr = "E2"; Val = "0.93%" },
    @{ Addr = "D3"; Val = "36.40" },
    @{ Addr = "E3"; Val = "1.61%" },
    @{ Addr = "D4"; Val = "5.062" },
    @{ Addr = "E4"; Val = "1.58%" },
    @{ Addr = "D5"; Val = "0.08081" },
    @{ Addr = "E5"; Val = "0.04%" },
    @{ Addr = "D6"; Val = "1.971" },
    @{ Addr = "E6"; Val = "3.95%" },
    @{ Addr = "D7"; Val = "7.858" },
    @{ Addr = "E7"; Val = "-0.17%" },
    @{ Addr = "D8"; Val = "0.9285" },
    @{ Addr = "E8"; Val = "-0.17%" },
    @{ Addr = "D9"; Val = "0.1456" },
    @{ Addr = "E9"; Val = "12.25%" },
    @{ Addr = "D10"; Val = "0.1931" },
    @{ Addr = "E10"; Val = "1.50%" },
    @{ Addr = "D11"; Val = "0.09102" },
    @{ Addr = "E11"; Val = "-0.94%" },
    @{ Addr = "D12"; Val = "0.03528" },
    @{ Addr = "E12"; Val = "0.63%" },
    @{ Addr = "D13"; Val = "0.09869" },
    @{ Addr = "E13"; Val = "-0.55%" },
    @{ Addr = "D14"; Val = "0.001417" },
    @{ Addr = "E14"; Val = "0.09%" },
    @{ Addr = "D15"; Val = "0.006360" },
    @{ Addr = "E15"; Val = "-2.08%" },
    @{ Addr = "E16"; Val = "6.27%" },
    @{ Addr = "D17"; Val = "4.165" },
    @{ Addr = "E17"; Val = "0.37%" },
    @{ Addr = "D18"; Val = "3.435" },
    @{ Addr = "E18"; Val = "6.78%" },
    @{ Addr = "E19"; Val = "0.00%" },
    @{ Addr = "D20"; Val = "0.1328" },
    @{ Addr = "E20"; Val = "2.43%" },
    @{ Addr = "D21"; Val = "4.803" },
    @{ Addr = "E21"; Val = "-8.68%" },
    @{ Addr = "D22"; Val = "0.2349" },
    @{ Addr = "D23"; Val = "0.04397" },
    @{ Addr = "E23"; Val = "-0.58%" },
    @{ Addr = "D24"; Val = "0.001237" },
    @{ Addr = "E24"; Val = "0.23%" },
    @{ Addr = "D25"; Val = "0.004166" },
    @{ Addr = "E25"; Val = "-11.73%" },
    @{ Addr = "D27"; Val = "0.0001303" },
    @{ Addr = "E27"; Val = "0.16%" },
    @{ Addr = "D39"; Val = "0.02041" },
    @{ Addr = "E39"; Val = "4.13%" },
    @{ Addr = "D40"; Val = "0.05114" },
    @{ Addr = "E40"; Val = "-2.18%" },
    @{ Addr = "D41"; Val = "0.007470" },
    @{ Addr = "E41"; Val = "-1.01%" },
    @{ Addr = "D42"; Val = "0.01002" },
    @{ Addr = "E42"; Val = "-1.73%" },
    @{ Addr = "D43"; Val = "0.1364" },
    @{ Addr = "E43"; Val = "-0.65%" },
    @{ Addr = "D44"; Val = "0.002125" },
    @{ Addr = "E44"; Val = "1.11%" },
    @{ Addr = "D45"; Val = "0.009883" },
    @{ Addr = "E45"; Val = "-7.75%" },
    @{ Addr = "D46"; Val = "0.00006304" },
    @{ Addr = "E46"; Val = "-0.58%" },
    @{ Addr = "E47"; Val = "0.17%" },
    @{ Addr = "D48"; Val = "64.80" },
    @{ Addr = "E48"; Val = "-0.24%" },
    @{ Addr = "D49"; Val = "0.001604" },
    @{ Addr = "E49"; Val = "-3.39%" },
    @{ Addr = "D50"; Val = "0.00002105" },
    @{ Addr = "E50"; Val = "0.17%" },
    @{ Addr = "D51"; Val = "0.0002005" },
    @{ Addr = "E51"; Val = "0.17%" }
)

$scratch = $ws.Range("ZZ1")
foreach ($u in $updates) {
    $scratch.Value = "'" + $u.Val
    $scratch.Copy()
    $ws.Range($u.Addr).PasteSpecial(-4163)  # xlPasteValues
}
$scratch.Clear()
